$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NutritionalData")
$ws.Range("A1").Value = "TEST"
